# The "Förändrad" (Changed) column (C) date stamp was bumped by one day
# (45177 -> 45178, i.e. 2023-09-08 -> 2023-09-09) for every data row
# (rows 2 through 267) on the only worksheet in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C267").Value = 45178
